$d = $word.ActiveDocument

# The target phrase spans three runs in the source document:
#   Run A (bold=0): " Classic & Quantum Mechanics"
#   Run B (plain) :  " "
#   Run C (bold=0, unchanged): "Engineer"
# We locate the whole phrase once, then rewrite Run A and Run B in place
# (via sub-ranges that line up exactly with the original run boundaries)
# so each keeps its own formatting and Run C is left untouched.

$oldA = " Classic & Quantum Mechanics"
$oldB = " "
$oldC = "Engineer"
$newA = " Material"
$newB = " Scientists & "

$whole = $d.Content
$found = $whole.Find.Execute($oldA + $oldB + $oldC, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if ($found) {
    $start = $whole.Start

    $rangeA = $d.Range($start, $start + $oldA.Length)
    $rangeA.Text = $newA

    $bStart = $rangeA.End
    $rangeB = $d.Range($bStart, $bStart + $oldB.Length)
    $rangeB.Text = $newB
}
